$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear all existing content (values) so the shared-string table rebuilds cleanly ---
$ws.Cells.ClearContents()

# --- Rewrite header row (row 1) to reseed shared-string slots 0-19 in original order ---
$ws.Range('A1').Value = 'Sending cluster'
$ws.Range('B1').Value = 'Ligand symbol'
$ws.Range('C1').Value = 'Receptor symbol'
$ws.Range('D1').Value = 'Target cluster'
$ws.Range('E1').Value = 'Ligand-expressing cells'
$ws.Range('F1').Value = 'Ligand detection rate'
$ws.Range('G1').Value = 'Ligand average expression value'
$ws.Range('H1').Value = 'Ligand total expression value'
$ws.Range('I1').Value = 'Ligand derived specificity of average expression value'
$ws.Range('J1').Value = 'Ligand derived specificity of total expression value'
$ws.Range('K1').Value = 'Receptor-expressing cells'
$ws.Range('L1').Value = 'Receptor detection rate'
$ws.Range('M1').Value = 'Receptor average expression value'
$ws.Range('N1').Value = 'Receptor total expression value'
$ws.Range('O1').Value = 'Receptor derived specificity of average expression value'
$ws.Range('P1').Value = 'Receptor derived specificity of total expression value'
$ws.Range('Q1').Value = 'Edge average expression weight'
$ws.Range('R1').Value = 'Edge total expression weight'
$ws.Range('S1').Value = 'Edge average expression derived specificity'
$ws.Range('T1').Value = 'Edge total expression derived specificity'

# --- Write the full data grid (rows 2-10) ---
# Row 2
$ws.Range('A2').Value = 'ECs'
$ws.Range('B2').Value = 'Ptn'
$ws.Range('C2').Value = 'Plxnb2'
$ws.Range('D2').Value = 'ECs'
$ws.Range('E2').Value = 1
$ws.Range('F2').Value = 0.3333333333333333
$ws.Range('G2').Value = 0.3206776666666667
$ws.Range('H2').Value = 0.962033
$ws.Range('I2').Value = 0.03001977461414601
$ws.Range('J2').Value = 0.03001977461414601
$ws.Range('K2').Value = 3
$ws.Range('L2').Value = 1
$ws.Range('M2').Value = 8.806900666666666
$ws.Range('N2').Value = 26.420702
$ws.Range('O2').Value = 0.1733678197953833
$ws.Range('P2').Value = 0.1733678197953834
$ws.Range('Q2').Value = 2.824176356351778
$ws.Range('R2').Value = 25.417587207166
$ws.Range('S2').Value = 0.00520446287560329
$ws.Range('T2').Value = 0.00520446287560329
# Row 3
$ws.Range('A3').Value = 'ECs'
$ws.Range('B3').Value = 'Ptn'
$ws.Range('C3').Value = 'Plxnb2'
$ws.Range('D3').Value = 'FAPs'
$ws.Range('E3').Value = 1
$ws.Range('F3').Value = 0.3333333333333333
$ws.Range('G3').Value = 0.3206776666666667
$ws.Range('H3').Value = 0.962033
$ws.Range('I3').Value = 0.03001977461414601
$ws.Range('J3').Value = 0.03001977461414601
$ws.Range('K3').Value = 3
$ws.Range('L3').Value = 1
$ws.Range('M3').Value = 18.76689066666667
$ws.Range('N3').Value = 56.30067200000001
$ws.Range('O3').Value = 0.3694347242421866
$ws.Range('P3').Value = 0.3694347242421866
$ws.Range('Q3').Value = 6.018122709575112
$ws.Range('R3').Value = 54.16310438617601
$ws.Range('S3').Value = 0.01109034715638963
$ws.Range('T3').Value = 0.01109034715638963
# Row 4
$ws.Range('A4').Value = 'ECs'
$ws.Range('B4').Value = 'Ptn'
$ws.Range('C4').Value = 'Plxnb2'
$ws.Range('D4').Value = 'sCs'
$ws.Range('E4').Value = 1
$ws.Range('F4').Value = 0.3333333333333333
$ws.Range('G4').Value = 0.3206776666666667
$ws.Range('H4').Value = 0.962033
$ws.Range('I4').Value = 0.03001977461414601
$ws.Range('J4').Value = 0.03001977461414601
$ws.Range('K4').Value = 3
$ws.Range('L4').Value = 1
$ws.Range('M4').Value = 23.225144
$ws.Range('N4').Value = 69.675432
$ws.Range('O4').Value = 0.4571974559624301
$ws.Range('P4').Value = 0.4571974559624301
$ws.Range('Q4').Value = 7.447784985917334
$ws.Range('R4').Value = 67.03006487325601
$ws.Range('S4').Value = 0.0137249645821531
$ws.Range('T4').Value = 0.0137249645821531
# Row 5
$ws.Range('A5').Value = 'FAPs'
$ws.Range('B5').Value = 'Ptn'
$ws.Range('C5').Value = 'Plxnb2'
$ws.Range('D5').Value = 'ECs'
$ws.Range('E5').Value = 3
$ws.Range('F5').Value = 1
$ws.Range('G5').Value = 5.752274333333333
$ws.Range('H5').Value = 17.256823
$ws.Range('I5').Value = 0.5384908178993973
$ws.Range('J5').Value = 0.5384908178993975
$ws.Range('K5').Value = 3
$ws.Range('L5').Value = 1
$ws.Range('M5').Value = 8.806900666666666
$ws.Range('N5').Value = 26.420702
$ws.Range('O5').Value = 0.1733678197953833
$ws.Range('P5').Value = 0.1733678197953834
$ws.Range('Q5').Value = 50.65970866108288
$ws.Range('R5').Value = 455.937377949746
$ws.Range('S5').Value = 0.09335697907905131
$ws.Range('T5').Value = 0.09335697907905134
# Row 6
$ws.Range('A6').Value = 'FAPs'
$ws.Range('B6').Value = 'Ptn'
$ws.Range('C6').Value = 'Plxnb2'
$ws.Range('D6').Value = 'FAPs'
$ws.Range('E6').Value = 3
$ws.Range('F6').Value = 1
$ws.Range('G6').Value = 5.752274333333333
$ws.Range('H6').Value = 17.256823
$ws.Range('I6').Value = 0.5384908178993973
$ws.Range('J6').Value = 0.5384908178993975
$ws.Range('K6').Value = 3
$ws.Range('L6').Value = 1
$ws.Range('M6').Value = 18.76689066666667
$ws.Range('N6').Value = 56.30067200000001
$ws.Range('O6').Value = 0.3694347242421866
$ws.Range('P6').Value = 0.3694347242421866
$ws.Range('Q6').Value = 107.9523034983396
$ws.Range('R6').Value = 971.5707314850562
$ws.Range('S6').Value = 0.1989372068176134
$ws.Range('T6').Value = 0.1989372068176134
# Row 7
$ws.Range('A7').Value = 'FAPs'
$ws.Range('B7').Value = 'Ptn'
$ws.Range('C7').Value = 'Plxnb2'
$ws.Range('D7').Value = 'sCs'
$ws.Range('E7').Value = 3
$ws.Range('F7').Value = 1
$ws.Range('G7').Value = 5.752274333333333
$ws.Range('H7').Value = 17.256823
$ws.Range('I7').Value = 0.5384908178993973
$ws.Range('J7').Value = 0.5384908178993975
$ws.Range('K7').Value = 3
$ws.Range('L7').Value = 1
$ws.Range('M7').Value = 23.225144
$ws.Range('N7').Value = 69.675432
$ws.Range('O7').Value = 0.4571974559624301
$ws.Range('P7').Value = 0.4571974559624301
$ws.Range('Q7').Value = 133.5973997191707
$ws.Range('R7').Value = 1202.376597472536
$ws.Range('S7').Value = 0.2461966320027327
$ws.Range('T7').Value = 0.2461966320027327
# Row 8
$ws.Range('A8').Value = 'sCs'
$ws.Range('B8').Value = 'Ptn'
$ws.Range('C8').Value = 'Plxnb2'
$ws.Range('D8').Value = 'ECs'
$ws.Range('E8').Value = 3
$ws.Range('F8').Value = 1
$ws.Range('G8').Value = 4.609262333333334
$ws.Range('H8').Value = 13.827787
$ws.Range('I8').Value = 0.4314894074864565
$ws.Range('J8').Value = 0.4314894074864565
$ws.Range('K8').Value = 3
$ws.Range('L8').Value = 1
$ws.Range('M8').Value = 8.806900666666666
$ws.Range('N8').Value = 26.420702
$ws.Range('O8').Value = 0.1733678197953833
$ws.Range('P8').Value = 0.1733678197953834
$ws.Range('Q8').Value = 40.59331551627488
$ws.Range('R8').Value = 365.339839646474
$ws.Range('S8').Value = 0.07480637784072873
$ws.Range('T8').Value = 0.07480637784072874
# Row 9
$ws.Range('A9').Value = 'sCs'
$ws.Range('B9').Value = 'Ptn'
$ws.Range('C9').Value = 'Plxnb2'
$ws.Range('D9').Value = 'FAPs'
$ws.Range('E9').Value = 3
$ws.Range('F9').Value = 1
$ws.Range('G9').Value = 4.609262333333334
$ws.Range('H9').Value = 13.827787
$ws.Range('I9').Value = 0.4314894074864565
$ws.Range('J9').Value = 0.4314894074864565
$ws.Range('K9').Value = 3
$ws.Range('L9').Value = 1
$ws.Range('M9').Value = 18.76689066666667
$ws.Range('N9').Value = 56.30067200000001
$ws.Range('O9').Value = 0.3694347242421866
$ws.Range('P9').Value = 0.3694347242421866
$ws.Range('Q9').Value = 86.50152226365157
$ws.Range('R9').Value = 778.5137003728642
$ws.Range('S9').Value = 0.1594071702681835
$ws.Range('T9').Value = 0.1594071702681835
# Row 10
$ws.Range('A10').Value = 'sCs'
$ws.Range('B10').Value = 'Ptn'
$ws.Range('C10').Value = 'Plxnb2'
$ws.Range('D10').Value = 'sCs'
$ws.Range('E10').Value = 3
$ws.Range('F10').Value = 1
$ws.Range('G10').Value = 4.609262333333334
$ws.Range('H10').Value = 13.827787
$ws.Range('I10').Value = 0.4314894074864565
$ws.Range('J10').Value = 0.4314894074864565
$ws.Range('K10').Value = 3
$ws.Range('L10').Value = 1
$ws.Range('M10').Value = 23.225144
$ws.Range('N10').Value = 69.675432
$ws.Range('O10').Value = 0.4571974559624301
$ws.Range('P10').Value = 0.4571974559624301
$ws.Range('Q10').Value = 107.0507814254427
$ws.Range('R10').Value = 963.4570328289841
$ws.Range('S10').Value = 0.1972758593775443
$ws.Range('T10').Value = 0.1972758593775443
